$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert a brand-new "Player Info" worksheet as the FIRST sheet.
# ---------------------------------------------------------------------------
$firstSheet = $wb.Worksheets.Item(1)
$playerInfo = $wb.Worksheets.Add($firstSheet)
$playerInfo.Name = "Player Info"

# Header row (bold, centered, bordered - matches the style used by the other
# sheets in this workbook for their header rows).
$playerInfo.Range("A1").Value = "ID"
$playerInfo.Range("B1").Value = "NAME"
$playerInfo.Range("C1").Value = "BATTING_HAND"
$playerInfo.Range("D1").Value = "BOWL_STYLE"

$headerRange = $playerInfo.Range("A1:D1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# Data row
$playerInfo.Range("A2").NumberFormat = "@"
$playerInfo.Range("A2").Value = "6033"
$playerInfo.Range("B2").Value = "Pothotuwa Arachchige Dhananjaya Lakshan"
$playerInfo.Range("C2").Value = "Left Handed"
$playerInfo.Range("D2").Value = "Right Arm Medium Fast"

$playerInfo.Range("A1").Select()

# ---------------------------------------------------------------------------
# 2. "ODI Batting" sheet: MATCH_CARD_LINK -> MATCH_CODE
#    (column D header text changes, and the URLs are replaced by just the
#    numeric match code that used to be the query-string parameter)
# ---------------------------------------------------------------------------
$batting = $wb.Worksheets.Item("ODI Batting")

$batting.Range("D1").Value = "MATCH_CODE"

$batting.Range("D2").NumberFormat = "@"
$batting.Range("D2").Value = "4469"

$batting.Range("D3").NumberFormat = "@"
$batting.Range("D3").Value = "4671"

$batting.Range("D4").NumberFormat = "@"
$batting.Range("D4").Value = "4674"

# ---------------------------------------------------------------------------
# 3. "ODI Bowling" sheet: MATCH_CARD_LINK -> MATCH_CODE
# ---------------------------------------------------------------------------
$bowling = $wb.Worksheets.Item("ODI Bowling")

$bowling.Range("B1").Value = "MATCH_CODE"

$bowling.Range("B2").NumberFormat = "@"
$bowling.Range("B2").Value = "4469"

$bowling.Range("B3").NumberFormat = "@"
$bowling.Range("B3").Value = "4671"

$bowling.Range("B4").NumberFormat = "@"
$bowling.Range("B4").Value = "4674"
